$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: Script 3
$ws.Range("A7").Value = 44918
$ws.Range("A2").Copy()
$ws.Range("A7").PasteSpecial(-4122)

$ws.Range("B7").Value = "Script 3"

$ws.Range("C7").Value = "12:22 to 1:22"
$ws.Range("C7").NumberFormat = "h:mm"

$ws.Range("D7").Value = "1h"

$ws.Range("F7").Value = 1
$ws.Range("G7").Formula = "=G6 + F7"

# Row 8: UI
$ws.Range("B8").Value = "UI"
$ws.Range("C8").Value = "1:25 to 2:43"
$ws.Range("D8").Value = "1h18m"
$ws.Range("F8").Value = 1.3
$ws.Range("G8").Formula = "=G7 + F8"

# Row 9: Second Ability
$ws.Range("B9").Value = "Second Ability"
$ws.Range("C9").Value = "2:51 to 5:18"
$ws.Range("D9").Value = "2h27m"
$ws.Range("F9").Value = 2.45
$ws.Range("G9").Formula = "=G8 + F9"

# Update selection
$ws.Range("G10").Select()
